$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 20329894
$ws.Range("I132").Value = 20918972
$ws.Range("J132").Value = 6750
$ws.Range("K132").Value = 62756916
$ws.Range("L132").Value = 20250
$ws.Range("M132").Value = -62754386
$ws.Range("N132").Value = -25310
$ws.Range("H137").Value = 341389.34
$ws.Range("I137").Value = 542080.6
$ws.Range("J137").Value = 1757.9231
$ws.Range("K137").Value = 1626241.8
$ws.Range("L137").Value = 5273.7693
$ws.Range("M137").Value = -1623691.8
$ws.Range("N137").Value = -10373.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3380.8857
$ws.Range("I32").Value = 3102.4922
$ws.Range("J32").Value = 7000
$ws.Range("K32").Value = 3102.4922
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = -2815.4922
$ws.Range("N32").Value = -7574
$ws.Range("H45").Value = 2245783.5
$ws.Range("I45").Value = 7857542.5
$ws.Range("J45").Value = 1080
$ws.Range("K45").Value = 7857542.5
$ws.Range("L45").Value = 1080
$ws.Range("M45").Value = -7857165.5
$ws.Range("N45").Value = -1834
$ws.Range("H74").Value = 5853.409
$ws.Range("I74").Value = 980.8
$ws.Range("K74").Value = 980.8
$ws.Range("M74").Value = -106.8
$ws.Range("H77").Value = 5853.409
$ws.Range("I77").Value = 980.8
$ws.Range("K77").Value = 4904
$ws.Range("M77").Value = -536
$ws.Range("H110").Value = 820
$ws.Range("I110").Value = 775.6667
$ws.Range("J110").Value = 877
$ws.Range("K110").Value = 775.6667
$ws.Range("L110").Value = 877
$ws.Range("M110").Value = 1269.3333
$ws.Range("N110").Value = -4967
$ws.Range("H122").Value = 85715420
$ws.Range("I122").Value = 85715420
$ws.Range("K122").Value = 257146260
$ws.Range("M122").Value = -257143810
$ws.Range("H132").Value = 4168251.8
$ws.Range("I132").Value = 6251057.5
$ws.Range("J132").Value = 2639.6
$ws.Range("K132").Value = 18753172.5
$ws.Range("L132").Value = 7918.799999999999
$ws.Range("M132").Value = -18750642.5
$ws.Range("N132").Value = -12978.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1009.13043
$ws.Range("I105").Value = 741.625
$ws.Range("J105").Value = 1620.5714
$ws.Range("K105").Value = 741.625
$ws.Range("L105").Value = 1620.5714
$ws.Range("M105").Value = 1005.375
$ws.Range("N105").Value = -5114.5714
$ws.Range("H107").Value = 47621010
$ws.Range("I107").Value = 90911040
$ws.Range("J107").Value = 1970
$ws.Range("K107").Value = 90911040
$ws.Range("L107").Value = 1970
$ws.Range("M107").Value = -90909120
$ws.Range("N107").Value = -5810
$ws.Range("H134").Value = 41731920
$ws.Range("I134").Value = 41731920
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 125195760
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("M134").Value = -125193225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 977.7778
$ws.Range("J16").Value = 1066.6666
$ws.Range("L16").Value = 1066.6666
$ws.Range("N16").Value = -1640.6666
$ws.Range("H31").Value = 10125.544
$ws.Range("I31").Value = 1042.7307
$ws.Range("J31").Value = 21933.2
$ws.Range("K31").Value = 1042.7307
$ws.Range("L31").Value = 21933.2
$ws.Range("M31").Value = -747.7307000000001
$ws.Range("N31").Value = -22523.2
$ws.Range("H34").Value = 10125.544
$ws.Range("I34").Value = 1042.7307
$ws.Range("J34").Value = 21933.2
$ws.Range("K34").Value = 1042.7307
$ws.Range("L34").Value = 21933.2
$ws.Range("M34").Value = -840.7307000000001
$ws.Range("N34").Value = -22337.2
$ws.Range("H107").Value = 280.8889
$ws.Range("I107").Value = 225.04347
$ws.Range("J107").Value = 379.69232
$ws.Range("K107").Value = 225.04347
$ws.Range("L107").Value = 379.69232
$ws.Range("M107").Value = 1694.95653
$ws.Range("N107").Value = -4219.69232
$ws.Range("H113").Value = 977.7778
$ws.Range("J113").Value = 1066.6666
$ws.Range("L113").Value = 1066.6666
$ws.Range("N113").Value = -5406.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2200.5
$ws.Range("I109").Value = 1467.6666
$ws.Range("J109").Value = 2933.3333
$ws.Range("K109").Value = 4402.9998
$ws.Range("L109").Value = 8799.999899999999
$ws.Range("M109").Value = -3362.9998
$ws.Range("N109").Value = -10879.9999
$ws.Range("H121").Value = 537.8
$ws.Range("I121").Value = 404.25
$ws.Range("J121").Value = 690.4286
$ws.Range("K121").Value = 1212.75
$ws.Range("L121").Value = 2071.2858
$ws.Range("M121").Value = 97.25
$ws.Range("N121").Value = -4691.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2040.5
$ws.Range("I113").Value = 1791.2
$ws.Range("J113").Value = 2289.8
$ws.Range("K113").Value = 1791.2
$ws.Range("L113").Value = 2289.8
$ws.Range("M113").Value = 378.8
$ws.Range("N113").Value = -6629.8
$ws.Range("H122").Value = 58824904
$ws.Range("I122").Value = 125000960
$ws.Range("K122").Value = 375002880
$ws.Range("M122").Value = -375000430
$ws.Range("H132").Value = 29415416
$ws.Range("I132").Value = 50000756
$ws.Range("J132").Value = 7787.5713
$ws.Range("K132").Value = 150002268
$ws.Range("L132").Value = 23362.7139
$ws.Range("M132").Value = -149999738
$ws.Range("N132").Value = -28422.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2198.3684
$ws.Range("I61").Value = 1488.6666
$ws.Range("J61").Value = 3415
$ws.Range("K61").Value = 1488.6666
$ws.Range("L61").Value = 3415
$ws.Range("M61").Value = -1286.6666
$ws.Range("N61").Value = -3819
$ws.Range("H113").Value = 2198.3684
$ws.Range("I113").Value = 1488.6666
$ws.Range("J113").Value = 3415
$ws.Range("K113").Value = 1488.6666
$ws.Range("L113").Value = 3415
$ws.Range("M113").Value = 681.3334
$ws.Range("N113").Value = -7755
$ws.Range("H132").Value = 10002560
$ws.Range("I132").Value = 13335213
$ws.Range("J132").Value = 4599.6
$ws.Range("K132").Value = 40005639
$ws.Range("L132").Value = 13798.8
$ws.Range("M132").Value = -40003109
$ws.Range("N132").Value = -18858.8
$ws.Range("H136").Value = 3671.8823
$ws.Range("I136").Value = 5404.241
$ws.Range("J136").Value = 1388.3182
$ws.Range("K136").Value = 16212.723
$ws.Range("L136").Value = 4164.9546
$ws.Range("M136").Value = -13662.723
$ws.Range("N136").Value = -9264.954600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 615.3333
$ws.Range("I113").Value = 524.5263
$ws.Range("J113").Value = 831
$ws.Range("K113").Value = 1573.5789
$ws.Range("L113").Value = 2493
$ws.Range("M113").Value = 596.4211
$ws.Range("N113").Value = -6833
$ws.Range("H122").Value = 13434.25
$ws.Range("I122").Value = 25673.5
$ws.Range("J122").Value = 1195
$ws.Range("K122").Value = 77020.5
$ws.Range("L122").Value = 3585
$ws.Range("M122").Value = -74570.5
$ws.Range("N122").Value = -8485
$ws.Range("H132").Value = 51944080
$ws.Range("I132").Value = 50003480
$ws.Range("J132").Value = 57119012
$ws.Range("K132").Value = 150010440
$ws.Range("L132").Value = 171357036
$ws.Range("M132").Value = -150007910
$ws.Range("N132").Value = -171362096
$ws.Range("H136").Value = 18236930
$ws.Range("I136").Value = 13734717
$ws.Range("J136").Value = 23811098
$ws.Range("K136").Value = 41204151
$ws.Range("L136").Value = 71433294
$ws.Range("M136").Value = -41201601
$ws.Range("N136").Value = -71438394
